# Auto-generated: apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.632.55'
$ws.Range('E2').Value = '  -6.01%  '

$ws.Range('D3').Value = '2.217.24'
$ws.Range('E3').Value = '  -6.33%  '

$ws.Range('E4').Value = '  +0.21%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -11.04%  '

$ws.Range('E7').Value = '  -9.21%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.560'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -9.27%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.41'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -12.32%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.72'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.28%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0820'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -10.78%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -10.10%  '

$ws.Range('E14').Value = '  -4.08%  '

$ws.Range('D15').Value = '2.554.80'
$ws.Range('E15').Value = '  -6.28%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.854'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -13.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.01'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -9.32%  '

$ws.Range('D18').Value = '2.185.97'
$ws.Range('E18').Value = '  -7.62%  '

$ws.Range('D19').Value = '42.545.86'
$ws.Range('E19').Value = '  -6.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.62'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.28%  '

$ws.Range('D21').Value = '0.0₃0962'
$ws.Range('E21').Value = '  -10.05%  '

$ws.Range('E22').Value = '  -11.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '65.41'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -11.15%  '

$ws.Range('E24').Value = '  -7.69%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '235.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.65%  '

$ws.Range('E26').Value = '  -7.20%  '

$ws.Range('E27').Value = '  -0.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.01%  '

$ws.Range('E29').Value = '  -5.11%  '

$ws.Range('E30').Value = '  -13.43%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.45'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -9.01%  '

$ws.Range('E32').Value = '  -9.13%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.88'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -7.24%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '33.70'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -11.14%  '

$ws.Range('E35').Value = '  -7.52%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.22'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.80%  '

$ws.Range('E37').Value = '  -7.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.42'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -8.35%  '

$ws.Range('E40').Value = '  -12.09%  '

$ws.Range('E41').Value = '  -10.41%  '

$ws.Range('E42').Value = '  -10.74%  '

$ws.Range('E43').Value = '  +0.36%  '

$ws.Range('D44').Value = '1.785.08'
$ws.Range('E44').Value = '  +9.61%  '

$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '88.10'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -12.83%  '

$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.92'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.01%  '

$ws.Range('E47').Value = '  -11.92%  '

$ws.Range('E48').Value = '  -5.63%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.33'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '60.07'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -14.34%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '15.57'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +54.68%  '

